$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the "availability" sentence into two runs and change the
#    wording from "...internship objective." to
#    "...internship opportunities".
# ------------------------------------------------------------------
$full = $d.Content
$found = $full.Find.Execute(" for full-time, part-time  and internship objective.")

$matchStart = $full.Start
$matchEnd = $full.End

$splitPoint = $matchStart + " for full-time, part".Length

$r1 = $d.Range($matchStart, $splitPoint)
$r2 = $d.Range($splitPoint, $matchEnd)

# Clear the trailing portion of the original run, then append the
# replacement text as a brand-new run right after $r1.
$r2.Text = ""
$r1.InsertAfter("-time  and internship opportunities")

# ------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the paragraph that used
#    to hold the "Career objective" heading text up onto that same
#    paragraph (which becomes empty), leaving the paragraph that used
#    to hold the bookmark empty too.
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$headingPara = $paras.Item(30)   # "Career objective" paragraph

# Re-anchor the bookmark at the very start of the heading paragraph
# *before* touching the text, so the engine doesn't snap it back to
# the top of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bookmarkTarget = $d.Range($headingPara.Range.Start, $headingPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkTarget)

# Now remove the "Career objective" run text, keeping the paragraph
# (and its underline formatting) intact but empty.
$textRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$textRange.Delete()
